{"js": "// Convert the field (fldChar begin/instrText.../fldChar end) that builds\n// the text:\n//   m:'https://www.m2doc.org/tests/'.fromHTMLURI()\n// into literal text runs wrapped in \"{\" ... \"}\" (no more field codes),\n// matching the new TokenIteratorFieldRewriterSplit tokenizer output.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that holds the field (search instead of hard-coding\n// an index so the script stays correct even if the document layout shifts).\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].fields.load(\"items\");\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].fields.items.length > 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Replace only the paragraph's inline content (keeps the paragraph mark\n  // itself untouched) with plain w:t runs instead of fldChar/instrText\n  // runs, preserving the existing bookmarkStart/bookmarkEnd pair.\n  const contentRange = target.getRange(\"Content\");\n\n  const runsXml =\n    \"<w:r><w:t>{</w:t></w:r>\" +\n    \"<w:r><w:t>m</w:t></w:r>\" +\n    \"<w:r><w:t>:</w:t></w:r>\" +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    \"<w:r><w:t>http</w:t></w:r>\" +\n    \"<w:r><w:t>s</w:t></w:r>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>://www.m2doc.org/tests/</w:t></w:r>\" +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    \"<w:r><w:t>.fromHTML</w:t></w:r>\" +\n    \"<w:r><w:t>URI</w:t></w:r>\" +\n    \"<w:r><w:t>()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>';\n\n  const ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" + runsXml + \"</w:p></w:body>\" +\n    \"</w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\n  contentRange.insertOoxml(ooxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Convert the field (fldChar begin/instrText.../fldChar end) that builds\n# the text:\n#   m:'https://www.m2doc.org/tests/'.fromHTMLURI()\n# into literal text runs wrapped in \"{\" ... \"}\" (no more field codes),\n# matching the new TokenIteratorFieldRewriterSplit tokenizer output.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that holds the field (search instead of hard-coding\n# an index so the script stays correct even if the document layout shifts).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:r><w:t>{</w:t></w:r>' +\n        '<w:r><w:t>m</w:t></w:r>' +\n        '<w:r><w:t>:</w:t></w:r>' +\n        \"<w:r><w:t>'</w:t></w:r>\" +\n        '<w:r><w:t>http</w:t></w:r>' +\n        '<w:r><w:t>s</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n        '<w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>://www.m2doc.org/tests/</w:t></w:r>' +\n        \"<w:r><w:t>'</w:t></w:r>\" +\n        '<w:r><w:t>.fromHTML</w:t></w:r>' +\n        '<w:r><w:t>URI</w:t></w:r>' +\n        '<w:r><w:t>()</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n        '</w:p>'\n\n    $target.Range.InsertXML($xml)\n}\n"}
